$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 ("I0") and J1 ("IF")
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy formatting (bold font, thin border, centered alignment) from H1 to the new headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Fill in data values for I2:J46
$data = New-Object "object[,]" 45,2
$data[0,0] = 8
$data[0,1] = 8
$data[1,0] = 8
$data[1,1] = 8
$data[2,0] = 7
$data[2,1] = 7
$data[3,0] = 8
$data[3,1] = 8
$data[4,0] = 8
$data[4,1] = 8
$data[5,0] = 8
$data[5,1] = 8
$data[6,0] = 8
$data[6,1] = 8
$data[7,0] = 7
$data[7,1] = 7
$data[8,0] = 7
$data[8,1] = 7
$data[9,0] = 8
$data[9,1] = 9
$data[10,0] = 7
$data[10,1] = 7
$data[11,0] = 8
$data[11,1] = 8
$data[12,0] = 7
$data[12,1] = 7
$data[13,0] = 5
$data[13,1] = 7
$data[14,0] = 6
$data[14,1] = 7
$data[15,0] = 9
$data[15,1] = 9
$data[16,0] = 8
$data[16,1] = 8
$data[17,0] = 6
$data[17,1] = 6
$data[18,0] = 6
$data[18,1] = 6
$data[19,0] = 8
$data[19,1] = 8
$data[20,0] = 7
$data[20,1] = 7
$data[21,0] = 7
$data[21,1] = 8
$data[22,0] = 7
$data[22,1] = 8
$data[23,0] = 8
$data[23,1] = 8
$data[24,0] = 8
$data[24,1] = 8
$data[25,0] = 3
$data[25,1] = 6
$data[26,0] = 9
$data[26,1] = 9
$data[27,0] = 6
$data[27,1] = 8
$data[28,0] = 9
$data[28,1] = 9
$data[29,0] = 9
$data[29,1] = 9
$data[30,0] = 8
$data[30,1] = 8
$data[31,0] = 8
$data[31,1] = 8
$data[32,0] = 8
$data[32,1] = 9
$data[33,0] = 7
$data[33,1] = 7
$data[34,0] = 9
$data[34,1] = 9
$data[35,0] = 5
$data[35,1] = 6
$data[36,0] = 4
$data[36,1] = 4
$data[37,0] = 6
$data[37,1] = 6
$data[38,0] = 2
$data[38,1] = 4
$data[39,0] = 6
$data[39,1] = 6
$data[40,0] = 7
$data[40,1] = 7
$data[41,0] = 5
$data[41,1] = 5
$data[42,0] = 5
$data[42,1] = 5
$data[43,0] = 8
$data[43,1] = 8
$data[44,0] = 5
$data[44,1] = 5

$ws.Range("I2:J46").Value = $data
